$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.722.61'
$ws.Range("E2").Value = '  +0.92%  '

# Row 3
$ws.Range("D3").Value = '1.618.20'
$ws.Range("E3").Value = '  +0.77%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.993'
$ws.Range("E4").Value = '  -0.61%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.66'
$ws.Range("E5").Value = '  +0.03%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.520'
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.992'
$ws.Range("E7").Value = '  -0.74%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.07'
$ws.Range("E8").Value = '  +8.24%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.258'
$ws.Range("E9").Value = '  +3.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0607'
$ws.Range("E10").Value = '  +1.10%  '

# Row 11
$ws.Range("E11").Value = '  -0.05%  '

# Row 12
$ws.Range("D12").Value = '1.853.86'
$ws.Range("E12").Value = '  +0.97%  '

# Row 13
$ws.Range("D13").Value = '1.636.77'
$ws.Range("E13").Value = '  +1.61%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.565'
$ws.Range("E14").Value = '  +5.55%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.89'
$ws.Range("E15").Value = '  +5.07%  '

# Row 16
$ws.Range("D16").Value = '29.766.43'
$ws.Range("E16").Value = '  +1.01%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.91'
$ws.Range("E17").Value = '  +16.63%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.30'
$ws.Range("E18").Value = '  +1.91%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.88'
$ws.Range("E19").Value = '  -0.36%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0707'
$ws.Range("E20").Value = '  +2.63%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.993'
$ws.Range("E21").Value = '  -0.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.08'
$ws.Range("E22").Value = '  +2.12%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.62'
$ws.Range("E23").Value = '  +4.50%  '

# Row 24
$ws.Range("E24").Value = '  +0.90%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.09'
$ws.Range("E25").Value = '  +0.30%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.62'
$ws.Range("E26").Value = '  +2.23%  '

# Row 27
$ws.Range("E27").Value = '  +1.99%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.59'
$ws.Range("E28").Value = '  +3.22%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.993'
$ws.Range("E29").Value = '  -0.63%  '

# Row 30
$ws.Range("E30").Value = '  +3.01%  '

# Row 31
$ws.Range("E31").Value = '  +2.30%  '

# Row 32
$ws.Range("E32").Value = '  +2.90%  '

# Row 34
$ws.Range("D34").Value = '1.422.83'
$ws.Range("E34").Value = '  +0.69%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.62'
$ws.Range("E35").Value = '  +6.00%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.03'
$ws.Range("E36").Value = '  -0.35%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.87'
$ws.Range("E37").Value = '  +1.59%  '

# Row 38
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.29'
$ws.Range("E38").Value = '  -0.68%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0170'
$ws.Range("E39").Value = '  +2.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.558'
$ws.Range("E40").Value = '  +3.80%  '

# Row 41
$ws.Range("E41").Value = '  +2.70%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.825'
$ws.Range("E42").Value = '  +3.35%  '

# Row 43
$ws.Range("E43").Value = '  -0.14%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '69.35'
$ws.Range("E44").Value = '  +5.63%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '53.38'
$ws.Range("E45").Value = '  +1.37%  '

# Row 46
$ws.Range("E46").Value = '  +18.55%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.991'
$ws.Range("E47").Value = '  -0.73%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.44'
$ws.Range("E48").Value = '  +2.95%  '

# Row 49
$ws.Range("D49").Value = '1.761.74'
$ws.Range("E49").Value = '  +0.83%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.10'
$ws.Range("E50").Value = '  +1.47%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0533'
$ws.Range("E51").Value = '  +2.00%  '
